# Update "想去人数" (F column, people interested) and a couple of
# "最低票价" (G column, min ticket price) figures across the four sheets
# of the 北京-漫展信息 workbook, matching the regenerated site export
# at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 748
$ws.Range("F6").Value = 2353
$ws.Range("F8").Value = 1750
$ws.Range("F9").Value = 2964
$ws.Range("F10").Value = 171
$ws.Range("F11").Value = 4423
$ws.Range("G11").Value = 70
$ws.Range("F12").Value = 385
$ws.Range("F13").Value = 210
$ws.Range("F17").Value = 17
$ws.Range("F18").Value = 207
$ws.Range("F19").Value = 83
$ws.Range("F21").Value = 305
$ws.Range("F22").Value = 4454
$ws.Range("F24").Value = 3672
$ws.Range("F25").Value = 1137
$ws.Range("F27").Value = 563
$ws.Range("F28").Value = 4373
$ws.Range("F29").Value = 90
$ws.Range("F30").Value = 575
$ws.Range("F31").Value = 575
$ws.Range("F32").Value = 538

# ---- 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 27
$ws.Range("F5").Value = 31

# ---- 本地生活 (Local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 209

# ---- 全部类型 (All types, combined) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 209
$ws.Range("F9").Value = 2353
$ws.Range("F11").Value = 1750
$ws.Range("F13").Value = 2964
$ws.Range("F14").Value = 171
$ws.Range("F15").Value = 4423
$ws.Range("G15").Value = 70
$ws.Range("F16").Value = 385
$ws.Range("F17").Value = 210
$ws.Range("F21").Value = 17
$ws.Range("F22").Value = 207
$ws.Range("F23").Value = 27
$ws.Range("F24").Value = 83
$ws.Range("F26").Value = 305
$ws.Range("F27").Value = 4454
$ws.Range("F29").Value = 3672
$ws.Range("F30").Value = 1137
$ws.Range("F32").Value = 563
$ws.Range("F33").Value = 4373
$ws.Range("F34").Value = 90
$ws.Range("F35").Value = 575
$ws.Range("F36").Value = 575
$ws.Range("F37").Value = 538
$ws.Range("F39").Value = 31
